# Add a new "volume_weighted_avg_spread_in_USD" column (F) to each of the
# four per-symbol pivot tables on the Pivots sheet, shifting the existing
# PnL_per_lot..pct_impact_on_PnL_exec_spread columns one slot to the right,
# and refresh every summary value to match the recomputed pivot output.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank column at F; this pushes old F:K (PnL_per_lot .. pct_impact_on_PnL_exec_spread)
# to G:L and extends the used range to column L.
$ws.Columns("F:F").Insert()

# Column widths: new F column -> 35 chars; J (old I, "total_volume") narrows 20 -> 19.
# (ColumnWidth assignments are rounded by the host to the nearest 1/6 character, so the
# inputs below are chosen so the stored width lands exactly on the target integer.)
$ws.Columns("F:F").ColumnWidth = 34.16666666666667   # -> stored width 35
$ws.Columns("J:J").ColumnWidth = 18.166666666666664  # -> stored width 19

# New column header text, one per pivot table (EURUSD, GBPUSD, USDJPY, XAUUSD).
$ws.Range("F2").Value = "volume_weighted_avg_spread_in_USD"
$ws.Range("F11").Value = "volume_weighted_avg_spread_in_USD"
$ws.Range("F20").Value = "volume_weighted_avg_spread_in_USD"
$ws.Range("F29").Value = "volume_weighted_avg_spread_in_USD"

# Refreshed pivot values (B:L) for every data row across the four tables.
function Set-RowValues($row, $vals) {
    $arr = New-Object 'object[,]' 1,11
    for ($i = 0; $i -lt $vals.Length; $i++) { $arr[0,$i] = $vals[$i] }
    $ws.Range("B$($row):L$($row)").Value = $arr
}

Set-RowValues 4 @(226, 26.71394799054373, 9.083333333331161, 8.800082670593104, 7.145760355215238, -7.246439295241086, -26774974.35, 23.97303753726283, 3694914.6, 28.9410265661055, 0.1379988100716892)
Set-RowValues 5 @(180, 21.27659574468085, 9.093750000000247, 8.767938246146342, 7.593880954716576, -23.39170510734847, -71235880.26000001, 63.78121633841791, 3045347.91, 23.85318858691453, 0.04275019693565867)
Set-RowValues 6 @(74, 8.747044917257684, 9.52631578947196, 9.379946294170733, 6.577107975926827, -7.069095899808215, -7893483.26, 7.067449179152464, 1116618.5, 8.746098129765944, 0.1414608054796838)
Set-RowValues 7 @(165, 19.50354609929078, 11.08695652173802, 9.446417343921874, 7.565519420988649, 1.520985848618986, 3384124.43, -3.029983929420951, 2224954.58, 17.4273228420917, 0.6574683130076278)
Set-RowValues 8 @(113, 13.35697399527187, 8.437500000000846, 8.534704789202021, 7.444681455572622, -5.534657292106162, -8592906.619999999, 7.693679562963027, 1552563.45, 12.16070868106516, 0.1806796604057592)
Set-RowValues 9 @(88, 10.40189125295508, 11.99999999999634, 10.43757810051786, 7.010237214750169, -0.5074365661726118, -574747.2299999999, 0.5146013116247049, 1132648.43, 8.871655194057174, 1.97068967170142)
Set-RowValues 13 @(194, 22.93144208037825, 9.53846153846129, 12.13233169369355, 10.16730385952625, -5.397910605132456, -9331144.373677, 9.728966156194963, 1728658.56, 22.10888035498022, 0.1852568656933995)
Set-RowValues 14 @(148, 17.49408983451537, 12.03225806451507, 12.57051141266667, 10.58948804115388, -41.03899665831257, -52182367.9508599, 54.40709857369549, 1271531.28, 16.26239766928574, 0.02436706745844497)
Set-RowValues 15 @(124, 14.65721040189125, 11.12765957446561, 11.9964327658, 9.790757268098226, -7.823823050719136, -10358295.1353324, 10.79990821831177, 1323942.92, 16.93272245451591, 0.1278147516268385)
Set-RowValues 16 @(201, 23.75886524822695, 12.55555555555693, 14.03384689476582, 11.27023592750013, 3.751553558781449, 6683405.68286, -6.968344405868774, 1781503.47, 22.78474649742982, 0.2665562371245507)
Set-RowValues 17 @(55, 6.501182033096926, 10.6666666666663, 15.26724701789474, 11.25885038668071, -44.44149006024582, -21707835.12738, 22.63332179001721, 488458.76, 6.247200304947436, 0.02250149575642894)
Set-RowValues 18 @(124, 14.65721040189125, 13.92307692307557, 13.25599805002062, 11.68507421398226, -7.360470106054805, -9014718.318076471, 9.399049667649333, 1224747.63, 15.66405271884089, 0.1358608873606305)
Set-RowValues 22 @(224, 26.47754137115839, 13.1851851851844, 10.34685395133333, 9.163127877642182, -13.18924566073865, -10835274.73, 24.12401358546816, 821523.46, 31.51593716196982, 0.07581934749890738)
Set-RowValues 23 @(187, 22.10401891252955, 13.88461538461494, 10.04526055128462, 9.149045741219078, -35.88564209826237, -27103837.67, 60.34487952197994, 755283.62, 28.97479167227323, 0.02786629809386694)
Set-RowValues 24 @(101, 11.93853427895981, 12.06896551723818, 10.55060927626582, 10.22009300772916, -9.084252122986966, -3171260, 7.060598021685028, 349094.23, 13.39223083938009, 0.1100806083386414)
Set-RowValues 25 @(166, 19.62174940898345, 13.81818181817928, 10.30664709236364, 8.111806541144686, -4.468655591029096, -1399487.57, 3.115865355762311, 313178.66, 12.01440914302059, 0.2237809514806909)
Set-RowValues 26 @(71, 8.392434988179669, 15.2500000000017, 12.00305448664, 8.342440709459956, -13.16094076206549, -1848576.79, 4.115732430140342, 140459.32, 5.388412283360735, 0.0759824102303048)
Set-RowValues 27 @(97, 11.46572104018913, 14.94999999999886, 10.633229949, 9.379553441179803, -2.449696561161502, -556455.58, 1.238911084964208, 227152.86, 8.714218899995537, 0.408213823644288)
Set-RowValues 31 @(234, 27.9236276849642, 20.76595744681236, 27.35697571637037, 24.5843936986288, -46.45713917658892, -166668754.47, 38.23268832562583, 3587581.1, 26.45322861051847, 0.02152521695748171)
Set-RowValues 32 @(161, 19.21241050119332, 21.12000000000353, 30.79435795365854, 29.11358610980235, -86.50893674627886, -273422756.32, 62.72133644220483, 3160630.18, 23.30508227530915, 0.01155949937210406)
Set-RowValues 33 @(59, 7.040572792362769, 20.6666666666706, 26.53536187878788, 23.96186508307103, -44.66863611574848, -39888229.5, 9.150090856466061, 892980.6900000001, 6.584442742590131, 0.02238707260747184)
Set-RowValues 34 @(216, 25.77565632458234, 21.94444444444722, 27.31399779281879, 24.87698885614075, 15.08552082243324, 47854004.64, -10.9773859554671, 3172181.14, 23.39025392078107, 0.06628872889247081)
Set-RowValues 35 @(121, 14.43914081145585, 22.17857142857399, 29.43266380088235, 27.73715359764924, -1.267282281380241, -2514695.109999999, 0.5768541001001538, 1984321.21, 14.63150271497155, 0.789090177218343)
Set-RowValues 36 @(47, 5.608591885441527, 21.50000000000318, 26.96942886, 25.67835588685672, -1.690700297241644, -1292175, 0.2964162310702217, 764283.89, 5.635489735829623, 0.5914708843616383)
